$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Component List")
$ws.Activate()

# --- Row 49: update the section header text ---
$ws.Range("B49").Value2 = "Hardware for ms41 case"

# --- Row 50: replace the placeholder row with the new "Thermal pad" component row ---
$ws.Rows(50).RowHeight = 16.5

$ws.Range("A50").Value2 = 14
$ws.Range("B50").Value2 = "Thermal pads"
$ws.Range("C50").Value2 = "Thermal pad"
$ws.Range("D50").Value2 = "Adhesive Thermal Pad for TO-220"
$ws.Range("E50").Value2 = ""
$ws.Range("F50").Value2 = ""
$ws.Range("G50").Value2 = "Aavid"
$ws.Range("H50").Value2 = "53-77-9ACG"
$ws.Range("I50").Value2 = "53-77-9ACG-ND"
$ws.Range("J50").Value2 = "532-53-77-9ACG"
$ws.Range("K50").Value2 = 0.55
$ws.Range("L50").Value2 = 0.673

$ws.Range("M50").Formula = "=K50*A50"
$ws.Range("N50").Formula = "=L50*A50"
$ws.Range("P50").Formula = '=IF(NOT(I50=""),A50&","&I50,"")'
$ws.Range("Q50").Formula = '=A50&"x "&C50'
$ws.Range("R50").Formula = '=IF(NOT(J50=""),J50&"|"&A50,"")'
$ws.Range("S50").Formula = '=H50&" "&A50'

# --- Update the view / selection state to match where the edit happened ---
$ws.Range("A50:XFD50").Select()

$wb.Save()
